# Updates cryptos list prices (D) and 1h volume % changes (E) to latest scraped values.
# Numeric-looking price strings are entered with a leading apostrophe so Excel
# stores them as text (matching the source data) instead of coercing to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.517.34'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.812.21'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D6').Value = '''305.55'
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('D7').Value = '''0.4549'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('D8').Value = '''0.3596'
$ws.Range('E8').Value = '  -1.82%  '
$ws.Range('D9').Value = '''46.28'
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('D10').Value = '''0.07109'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').Value = '''0.8921'
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('D12').Value = '''0.07715'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D14').Value = '1.812.86'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = '''5.258'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').Value = '''6.298'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '''85.99'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').Value = '''0.000008538'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').Value = '26.560.56'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '''14.12'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').Value = '''4.958'
$ws.Range('E23').Value = '  -1.10%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').Value = '''1.920'
$ws.Range('E25').Value = '  -3.17%  '
$ws.Range('D26').Value = '''152.04'
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').Value = '''17.79'
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('D29').Value = '''112.32'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').Value = '''4.819'
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').Value = '''0.08710'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').Value = '''3.132'
$ws.Range('E32').Value = '  +2.28%  '
$ws.Range('D33').Value = '''0.7386'
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('E34').Value = '  -1.82%  '
$ws.Range('D35').Value = '''2.716'
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('E38').Value = '  -0.96%  '
$ws.Range('D39').Value = '''2.914'
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('E40').Value = '  -0.59%  '
$ws.Range('D41').Value = '''0.5077'
$ws.Range('D42').Value = '''6.784'
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range('E43').Value = '  -3.45%  '
$ws.Range('D44').Value = '''8.018'
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('D45').Value = '''0.4675'
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('D47').Value = '''9.949'
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').Value = '''98.93'
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('D49').Value = '''1.565'
$ws.Range('D51').Value = '''63.68'
$ws.Range('E51').Value = '  -1.12%  '
